$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: the "Price" (D) column holds text values that look numeric
# (e.g. "305.99", "26.627.57"). Excel auto-converts a plain .Value
# assignment of a numeric-looking string into a real number, so for
# those cells we force a Text number format while writing, then restore
# the default "Normal" style so no formatting side effect is left behind.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = "26.627.57"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = "1.820.88"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.13%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "305.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "0.4678"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.53%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = "0.3598"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = "46.26"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "0.07133"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.64%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = "0.9011"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.95%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = "0.07812"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = "19.43"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "1.786.83"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.47%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "5.249"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = "6.334"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = "87.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("E18").Value = "  +0.05%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "0.000008566"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "1.008"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "26.667.70"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  -0.29%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = "5.014"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "10.56"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "1.932"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = "151.90"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "17.90"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = "1.978"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "113.69"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = "4.803"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.99%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "0.08797"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +2.96%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = "2.766"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.98%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "0.7314"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.91%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = "4.440"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "1.123"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value2 = "1.076"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = "0.01928"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.72%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = "2.917"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "0.05114"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "0.5064"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.68%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "6.811"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "0.1496"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "7.992"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = "0.4677"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = "1.009"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "9.965"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "98.97"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "1.559"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.82%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "0.06015"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = "63.78"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
